$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ D = 79.17607379688789; E = 77.20070648193359; F = 83.27875585724033; G = 76.39346719085641; H = 2319000000; I = "BABA" }
    3 = @{ D = 78.44480390100094; E = 74.39910125732422; F = 81.08495478617667; G = 72.3762535582847; H = 2319000000; I = "BABA" }
    4 = @{ D = 56.50685557919599; E = 79.61293792724609; F = 80.19225248846217; G = 55.27225273647926; H = 2319000000; I = "BABA" }
    5 = @{ D = 74.24716389304717; E = 63.65806198120117; F = 74.72201144742951; G = 62.05307492879728; H = 2319000000; I = "BABA" }
    6 = @{ D = 74.3041334321496; E = 73.06953430175781; F = 77.66605454103397; G = 71.85392588296074; H = 2319000000; I = "BABA" }
    7 = @{ D = 75.13987283656358; E = 78.33084869384766; F = 80.7240769640117; G = 73.77230968469955; H = 2319000000; I = "BABA" }
    8 = @{ D = 100.1453337216834; E = 96.57448577880859; F = 103.5167538318781; G = 94.01980393904527; H = 2319000000; I = "BABA" }
    9 = @{ D = 84.5228567809512; E = 96.21360015869141; F = 99.30960795269368; G = 83.64913910148488; H = 2319000000; I = "BABA" }
    10 = @{ D = 103.3743002994452; E = 109.6897735595703; F = 110.155122085298; G = 101.389441208159; H = 2319000000; I = "BABA" }
    11 = @{ D = 134.6192713438163; E = 147.1552429199219; F = 152.3215862340008; G = 132.4777044005996; H = 2319000000; I = "BABA" }
    12 = @{ D = 165.7882518690359; E = 175.5890960693359; F = 175.8075218469968; G = 160.0995739877707; H = 2319000000; I = "BABA" }
    13 = @{ D = 167.5252489114074; E = 194.01318359375; F = 195.8271045068636; G = 166.8614072427723; H = 2319000000; I = "BABA" }
    14 = @{ D = 173.6137342738819; E = 169.5585327148438; F = 174.3924911011387; G = 157.7728292097308; H = 2319000000; I = "BABA" }
    15 = @{ D = 172.5215926331687; E = 177.8113861083984; F = 188.3720040842142; G = 171.9517698517012; H = 2319000000; I = "BABA" }
    16 = @{ D = 157.5733937448925; E = 135.1226043701172; F = 157.6018834357928; G = 123.5173302309971; H = 2319000000; I = "BABA" }
    17 = @{ D = 127.3826060301531; E = 160.0141296386719; F = 161.1917422663955; G = 123.2989142778932; H = 2319000000; I = "BABA" }
    18 = @{ D = 175.7790450214942; E = 176.2349090576172; F = 180.2426086473627; G = 167.8680838898237; H = 2319000000; I = "BABA" }
    19 = @{ D = 167.0181172565888; E = 164.4017028808594; F = 170.8311421216834; G = 156.6996758113333; H = 2319000000; I = "BABA" }
    20 = @{ D = 159.5582606533422; E = 167.7826232910156; F = 170.6032187579979; G = 153.5466892283192; H = 2319000000; I = "BABA" }
    21 = @{ D = 205.7039511156663; E = 196.197494506836; F = 219.5125105979933; G = 189.4641601624099; H = 2319000000; I = "BABA" }
    22 = @{ D = 179.9672041436568; E = 192.4746856689453; F = 205.2290971708543; G = 175.7315580288462; H = 2319000000; I = "BABA" }
    23 = @{ D = 204.9536737729801; E = 238.3924407958984; F = 254.518257692142; G = 203.937507667616; H = 2319000000; I = "BABA" }
    24 = @{ D = 280.4069622200644; E = 289.3625793457031; F = 303.256622191907; G = 270.5491232281558; H = 2319000000; I = "BABA" }
    25 = @{ D = 215.1059151296817; E = 241.0610809326172; F = 255.46795218492; G = 209.0183667688753; H = 2319000000; I = "BABA" }
    26 = @{ D = 219.0756400548808; E = 219.3320617675781; F = 233.3305715662409; G = 210.1865085874636; H = 2319000000; I = "BABA" }
    27 = @{ D = 216.5969483061772; E = 185.3709716796875; F = 217.3377092561169; G = 170.6317010383922; H = 2319000000; I = "BABA" }
    28 = @{ D = 139.6906381535359; E = 156.6427001953125; F = 172.9299641247087; G = 131.466275735353; H = 2319000000; I = "BABA" }
    29 = @{ D = 113.3745868669021; E = 119.4621353149414; F = 131.7226940378181; G = 104.523429727425; H = 2319000000; I = "BABA" }
    30 = @{ D = 111.703131862575; E = 92.20588684082033; F = 112.9662207438479; G = 77.68505623092295; H = 2319000000; I = "BABA" }
    31 = @{ D = 107.9138498549603; E = 84.87425231933594; F = 119.5096259684844; G = 84.62732956857172; H = 2319000000; I = "BABA" }
    32 = @{ D = 75.74767652266307; E = 60.38161087036133; F = 81.94918381929942; G = 55.0918062970277; H = 2319000000; I = "BABA" }
    33 = @{ D = 86.52671304437573; E = 104.6563873291016; F = 115.198007749205; G = 86.18482226071843; H = 2319000000; I = "BABA" }
    34 = @{ D = 95.81472841424866; E = 80.42967224121094; F = 98.04651022546862; G = 78.07632378887907; H = 2319000000; I = "BABA" }
    35 = @{ D = 80.60061628474665; E = 97.02084350585938; F = 97.3437363240718; G = 78.48279327226724; H = 2319000000; I = "BABA" }
    36 = @{ D = 82.38604003370835; E = 78.38782501220703; F = 84.3291137143956; G = 75.04489460402499; H = 2319000000; I = "BABA" }
    37 = @{ D = 73.19519329466752; E = 69.46083831787109; F = 73.81116801069619; G = 64.12880121606212; H = 2319000000; I = "BABA" }
    38 = @{ D = 70.4906739185227; E = 72.04023742675781; F = 74.04215898737185; G = 65.79386472754183; H = 2319000000; I = "BABA" }
    39 = @{ D = 71.30481678480784; E = 77.53943634033203; F = 79.66353868268; G = 70.77379119922084; H = 2319000000; I = "BABA" }
    40 = @{ D = 104.6512005951864; E = 96.35148620605467; F = 115.8617236829262; G = 93.93236705074592; H = 2319000000; I = "BABA" }
    41 = @{ D = 82.98735743006301; E = 97.19718170166016; F = 101.9419901195062; G = 78.72932434290269; H = 2319000000; I = "BABA" }
    42 = @{ D = 130.5435703394336; E = 117.4449615478516; F = 131.8711380737685; G = 94.13888080467102; H = 2319000000; I = "BABA" }
    43 = @{ D = 112.9499969482422; E = 120.629997253418; F = 123.9899978637695; G = 103.7099990844727; H = 2319000000; I = "BABA" }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
    $ws.Range("H$row").Value = $vals.H
    $ws.Range("I$row").Value = $vals.I
}
